# Add a new "Resources" slide (slide 5) with a list of hyperlinked
# reference URLs, using the "Title and Content" layout (same layout
# used by the existing Flexbox / CSS Grid slides).

$p = $ppt.ActivePresentation

# Layout index 2 on the slide master == "Title and Content" (the same
# layout ppt/slideLayouts/slideLayout2.xml used by slides 2 and 3).
$s = $p.Slides.Add($p.Slides.Count + 1, 2)

# --- Title -----------------------------------------------------------
$title = $s.Shapes.Item(1)
$title.Name = "Title 6"
$title.TextFrame.TextRange.Text = "Resources"

# --- Body: one resource link per paragraph ----------------------------
$body = $s.Shapes.Item(2)
$body.Name = "Content Placeholder 7"

$tr = $body.TextFrame.TextRange

# Each entry is the full display text for one paragraph. The trailing
# "." is a placeholder for the last (blank) paragraph -- a genuinely
# empty trailing element would not be counted as its own paragraph by
# Paragraphs(), so it's typed then cleared below once the paragraph
# actually exists.
$lines = @(
    "https://github.com/shellwe/Flexbox_Grid",
    "https://css-tricks.com/snippets/css/a-guide-to-flexbox/",
    "http://flexboxfroggy.com/",
    "https://css-tricks.com/snippets/css/complete-guide-grid/",
    "http://cssgridgarden.com/",
    "https://tutorialzine.com/2017/03/css-grid-vs-flexbox",
    "."
)

$tr.Text = [string]::Join("`r", $lines)

# Clear the trailing placeholder paragraph so it stays a truly empty
# paragraph (no dangling run), matching a plain Enter press at the end.
$paraCount = $tr.Paragraphs().Count
$lastPara = $tr.Paragraphs($paraCount, 1)
$lastPara.Text = ""

# Run-split points (character counts) within each paragraph, matching
# how the links were typed/autocompleted in stages, and the matching
# hyperlink target for that whole paragraph.
$splits = @(
    @(8, 31),
    @(54, 1),
    @(24, 1),
    @(35, 20, 1),
    @(4, 20, 1),
    @(8, 44)
)
$urls = @(
    "https://github.com/shellwe/Flexbox_Grid",
    "https://css-tricks.com/snippets/css/a-guide-to-flexbox/",
    "http://flexboxfroggy.com/",
    "https://css-tricks.com/snippets/css/complete-guide-grid/",
    "http://cssgridgarden.com/",
    "https://tutorialzine.com/2017/03/css-grid-vs-flexbox"
)

for ($i = 0; $i -lt $urls.Length; $i++) {
    $para = $tr.Paragraphs($i + 1, 1)
    $url = $urls[$i]

    # First carve the paragraph into the target runs by touching each
    # sub-range individually (forces the run to split there).
    $pos = 1
    foreach ($len in $splits[$i]) {
        $chunk = $para.Characters($pos, $len)
        $chunk.ActionSettings.Item(1).Hyperlink.Address = $url
        $pos += $len
    }

    # Then apply the hyperlink to the whole paragraph too, so every run
    # (including ones already split above) ends up pointing at the URL.
    $whole = $para.Characters(1, $para.Length)
    $whole.ActionSettings.Item(1).Hyperlink.Address = $url
}
